$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.275.21'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '3.496.19'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.10'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.35'
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.34'
$ws.Range("E9").Value = '  +3.29%  '
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("E11").Value = '  +2.26%  '
$ws.Range("D12").Value = '4.096.00'
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("E14").Value = '  +0.83%  '
$ws.Range("D15").Value = '3.500.77'
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").Value = '64.374.32'
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.58'
$ws.Range("E17").Value = '  -6.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.84'
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.74'
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.49'
$ws.Range("E20").Value = '  -2.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '392.47'
$ws.Range("E21").Value = '  +2.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.571'
$ws.Range("E22").Value = '  +0.47%  '
$ws.Range("D23").Value = '3.638.51'
$ws.Range("E23").Value = '  -0.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.60'
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("E26").Value = '  +0.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.34'
$ws.Range("E28").Value = '  -2.03%  '
$ws.Range("E29").Value = '  +1.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.19'
$ws.Range("E30").Value = '  -2.16%  '
$ws.Range("E31").Value = '  -7.06%  '
$ws.Range("D32").Value = '3.518.65'
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("E33").Value = '  +5.41%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.40'
$ws.Range("E35").Value = '  -0.66%  '
$ws.Range("E36").Value = '  -4.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.86'
$ws.Range("E37").Value = '  -1.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.54'
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '166.81'
$ws.Range("E39").Value = '  +4.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0778'
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("E41").Value = '  -0.43%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.01'
$ws.Range("E43").Value = '  -6.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.38'
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("E45").Value = '  +3.24%  '
$ws.Range("E46").Value = '  -3.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.74'
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").Value = '2.374.27'
$ws.Range("E48").Value = '  -4.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.890'
$ws.Range("E49").Value = '  -2.80%  '
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.98'
$ws.Range("E51").Value = '  -1.38%  '
